$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'65.387.38"
$ws.Cells.Item(2, 5).Value = '  +0.07%  '
$ws.Cells.Item(3, 4).Value = "'3.549.37"
$ws.Cells.Item(3, 5).Value = '  +3.98%  '
$ws.Cells.Item(4, 4).Value = "'1.00"
$ws.Cells.Item(4, 5).Value = '  -0.17%  '
$ws.Cells.Item(5, 4).Value = "'600.02"
$ws.Cells.Item(5, 5).Value = '  +2.51%  '
$ws.Cells.Item(6, 4).Value = "'138.30"
$ws.Cells.Item(6, 5).Value = '  +0.74%  '
$ws.Cells.Item(7, 4).Value = "'3.549.04"
$ws.Cells.Item(7, 5).Value = '  +4.01%  '
$ws.Cells.Item(8, 5).Value = '  +0.02%  '
$ws.Cells.Item(9, 4).Value = "'0.493"
$ws.Cells.Item(9, 5).Value = '  -0.50%  '
$ws.Cells.Item(10, 4).Value = "'0.125"
$ws.Cells.Item(10, 5).Value = '  +3.71%  '
$ws.Cells.Item(11, 4).Value = "'6.89"
$ws.Cells.Item(11, 5).Value = '  -4.39%  '
$ws.Cells.Item(12, 4).Value = "'0.387"
$ws.Cells.Item(12, 5).Value = '  +3.42%  '
$ws.Cells.Item(13, 4).Value = "'4.152.55"
$ws.Cells.Item(13, 5).Value = '  +4.00%  '
$ws.Cells.Item(14, 4).Value = "'0.0000184"
$ws.Cells.Item(14, 5).Value = '  +2.90%  '
$ws.Cells.Item(15, 4).Value = "'27.27"
$ws.Cells.Item(15, 5).Value = '  +4.68%  '
$ws.Cells.Item(16, 4).Value = "'3.547.87"
$ws.Cells.Item(16, 5).Value = '  +3.64%  '
$ws.Cells.Item(17, 5).Value = '  +1.37%  '
$ws.Cells.Item(18, 4).Value = "'65.342.96"
$ws.Cells.Item(18, 5).Value = '  +0.00%  '
$ws.Cells.Item(19, 4).Value = "'10.25"
$ws.Cells.Item(19, 5).Value = '  +5.38%  '
$ws.Cells.Item(20, 4).Value = "'5.96"
$ws.Cells.Item(20, 5).Value = '  +1.66%  '
$ws.Cells.Item(21, 4).Value = "'14.33"
$ws.Cells.Item(21, 5).Value = '  +5.61%  '
$ws.Cells.Item(22, 4).Value = "'393.45"
$ws.Cells.Item(22, 5).Value = '  +1.55%  '
$ws.Cells.Item(23, 4).Value = "'0.575"
$ws.Cells.Item(23, 5).Value = '  +3.77%  '
$ws.Cells.Item(24, 4).Value = "'3.695.21"
$ws.Cells.Item(24, 5).Value = '  +3.96%  '
$ws.Cells.Item(25, 4).Value = "'73.81"
$ws.Cells.Item(25, 5).Value = '  +1.27%  '
$ws.Cells.Item(26, 4).Value = "'0.999"
$ws.Cells.Item(26, 5).Value = '  -0.06%  '
$ws.Cells.Item(27, 4).Value = "'0.0000116"
$ws.Cells.Item(27, 5).Value = '  +10.61%  '
$ws.Cells.Item(28, 4).Value = "'7.80"
$ws.Cells.Item(28, 5).Value = '  +10.29%  '
$ws.Cells.Item(29, 4).Value = "'1.00"
$ws.Cells.Item(29, 5).Value = '  -0.06%  '
$ws.Cells.Item(30, 4).Value = "'2.28"
$ws.Cells.Item(30, 5).Value = '  +3.36%  '
$ws.Cells.Item(31, 4).Value = "'8.16"
$ws.Cells.Item(31, 5).Value = '  +0.48%  '
$ws.Cells.Item(32, 4).Value = "'3.561.47"
$ws.Cells.Item(32, 5).Value = '  +4.16%  '
$ws.Cells.Item(33, 2).Value = 'USDe'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(33, 4).Value = "'1.00"
$ws.Cells.Item(33, 5).Value = '  +0.01%  '
$ws.Cells.Item(34, 2).Value = 'EthereumClassic'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(34, 4).Value = "'23.88"
$ws.Cells.Item(34, 5).Value = '  +4.61%  '
$ws.Cells.Item(35, 4).Value = "'0.146"
$ws.Cells.Item(35, 5).Value = '  +1.59%  '
$ws.Cells.Item(36, 4).Value = "'1.29"
$ws.Cells.Item(36, 5).Value = '  +13.01%  '
$ws.Cells.Item(37, 4).Value = "'6.97"
$ws.Cells.Item(37, 5).Value = '  +2.79%  '
$ws.Cells.Item(38, 2).Value = 'Monero'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(38, 4).Value = "'169.57"
$ws.Cells.Item(38, 5).Value = '  -1.08%  '
$ws.Cells.Item(39, 2).Value = 'ImmutableX'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(39, 4).Value = "'1.56"
$ws.Cells.Item(39, 5).Value = '  +7.61%  '
$ws.Cells.Item(40, 4).Value = "'5.01"
$ws.Cells.Item(40, 5).Value = '  +6.25%  '
$ws.Cells.Item(41, 4).Value = "'0.0804"
$ws.Cells.Item(41, 5).Value = '  +6.08%  '
$ws.Cells.Item(42, 4).Value = "'0.828"
$ws.Cells.Item(42, 5).Value = '  +1.62%  '
$ws.Cells.Item(43, 4).Value = "'26.59"
$ws.Cells.Item(43, 5).Value = '  +20.00%  '
$ws.Cells.Item(44, 4).Value = "'42.76"
$ws.Cells.Item(44, 5).Value = '  -1.67%  '
$ws.Cells.Item(45, 4).Value = "'1.00"
$ws.Cells.Item(45, 5).Value = '  -0.16%  '
$ws.Cells.Item(46, 4).Value = "'4.44"
$ws.Cells.Item(46, 5).Value = '  +1.25%  '
$ws.Cells.Item(47, 2).Value = 'ONDO'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Cells.Item(47, 4).Value = "'1.20"
$ws.Cells.Item(47, 5).Value = '  +10.07%  '
$ws.Cells.Item(48, 2).Value = 'Stacks'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(48, 4).Value = "'1.68"
$ws.Cells.Item(48, 5).Value = '  +5.28%  '
$ws.Cells.Item(49, 4).Value = "'6.80"
$ws.Cells.Item(49, 5).Value = '  +4.53%  '
$ws.Cells.Item(50, 4).Value = "'2.407.00"
$ws.Cells.Item(50, 5).Value = '  +10.36%  '
$ws.Cells.Item(51, 4).Value = "'307.92"
$ws.Cells.Item(51, 5).Value = '  +9.88%  '
